$d = $word.ActiveDocument
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()

$d.Content.Find.Execute("2024-04-02 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-03 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("38÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=", 2) | Out-Null
$d.Content.Find.Execute("39÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=", 2) | Out-Null
$d.Content.Find.Execute("79÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=", 2) | Out-Null
$d.Content.Find.Execute("98÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=", 2) | Out-Null
$d.Content.Find.Execute("91÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 2) | Out-Null
$d.Content.Find.Execute("15÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=", 2) | Out-Null
$d.Content.Find.Execute("87÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=", 2) | Out-Null
$d.Content.Find.Execute("39÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=", 2) | Out-Null
$d.Content.Find.Execute("14÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 2) | Out-Null
$d.Content.Find.Execute("23÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷4=", 2) | Out-Null
$d.Content.Find.Execute("33÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=", 2) | Out-Null
$d.Content.Find.Execute("94÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷7=", 2) | Out-Null
$d.Content.Find.Execute("50÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷4=", 2) | Out-Null
$d.Content.Find.Execute("46÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=", 2) | Out-Null
$d.Content.Find.Execute("80÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=", 2) | Out-Null
$d.Content.Find.Execute("37÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=", 2) | Out-Null
$d.Content.Find.Execute("86÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=", 2) | Out-Null
$d.Content.Find.Execute("69÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=", 2) | Out-Null
$d.Content.Find.Execute("72÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷6=", 2) | Out-Null
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=", 2) | Out-Null
$d.Content.Find.Execute("95÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=", 2) | Out-Null
$d.Content.Find.Execute("80÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=", 2) | Out-Null
$d.Content.Find.Execute("85÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=", 2) | Out-Null
$d.Content.Find.Execute("19÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷8=", 2) | Out-Null
$d.Content.Find.Execute("87÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=", 2) | Out-Null
